# Update FFXIV Leve profit calculations (scheduled market-price refresh).
$wb = $excel.ActiveWorkbook

# ALC row 21 (Leve Item ID 2149)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 2000
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2936

# ALC row 23 (Leve Item ID 2149)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2468

# ALC row 69 (Leve Item ID 12616)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9584.208000000001
$ws.Range("I69").Value = 4999.5
$ws.Range("J69").Value = 10001
$ws.Range("K69").Value = 14998.5
$ws.Range("L69").Value = 30003
$ws.Range("M69").Value = -14124.5
$ws.Range("N69").Value = -31751

# ALC row 72 (Leve Item ID 12616)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 9584.208000000001
$ws.Range("I72").Value = 4999.5
$ws.Range("J72").Value = 10001
$ws.Range("K72").Value = 44995.5
$ws.Range("L72").Value = 90009
$ws.Range("M72").Value = -40627.5
$ws.Range("N72").Value = -98745

# ALC row 116 (Leve Item ID 27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 379179.28
$ws.Range("I116").Value = 85644.14
$ws.Range("J116").Value = 892865.75
$ws.Range("K116").Value = 85644.14
$ws.Range("L116").Value = 892865.75
$ws.Range("M116").Value = -82202.14
$ws.Range("N116").Value = -899749.75

# ALC row 125 (Leve Item ID 36228)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2041.5
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 2222
$ws.Range("K125").Value = 13500
$ws.Range("L125").Value = 19998
$ws.Range("M125").Value = -11040
$ws.Range("N125").Value = -24918

# ALC row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 73205.39
$ws.Range("I132").Value = 79525.60000000001
$ws.Range("K132").Value = 238576.8
$ws.Range("M132").Value = -236046.8

# ALC row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2043.9333
$ws.Range("J138").Value = 2893.4783
$ws.Range("L138").Value = 8680.4349
$ws.Range("N138").Value = -18960.4349

# ARM row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19203.084
$ws.Range("I32").Value = 19203.084
$ws.Range("K32").Value = 19203.084
$ws.Range("M32").Value = -18916.084

# ARM row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5381.6
$ws.Range("I74").Value = 1960.7142
$ws.Range("K74").Value = 1960.7142
$ws.Range("M74").Value = -1086.7142

# ARM row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5381.6
$ws.Range("I77").Value = 1960.7142
$ws.Range("K77").Value = 9803.571
$ws.Range("M77").Value = -5435.571

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 446636.44
$ws.Range("I132").Value = 501841.2
$ws.Range("J132").Value = 4998.2
$ws.Range("K132").Value = 1505523.6
$ws.Range("L132").Value = 14994.6
$ws.Range("M132").Value = -1502993.6
$ws.Range("N132").Value = -20054.6

# BSM row 61 (Leve Item ID 2543)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# BSM row 132 (Leve Item ID 41855)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# CRP row 22 (Leve Item ID 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 683.3125
$ws.Range("I22").Value = 726.8570999999999
$ws.Range("J22").Value = 378.5
$ws.Range("K22").Value = 726.8570999999999
$ws.Range("L22").Value = 378.5
$ws.Range("M22").Value = -376.8570999999999
$ws.Range("N22").Value = -1078.5

# CRP row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 43704180
$ws.Range("I132").Value = 55557944
$ws.Range("K132").Value = 166673832
$ws.Range("M132").Value = -166671302

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4902.5
$ws.Range("I134").Value = 5185.5137
$ws.Range("K134").Value = 15556.5411
$ws.Range("M134").Value = -13021.5411

# CRP row 138 (Leve Item ID 42302)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

# CUL row 26 (Leve Item ID 4746)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 472.44446
$ws.Range("J26").Value = 446.66666
$ws.Range("L26").Value = 1339.99998
$ws.Range("N26").Value = -1915.99998

# CUL row 39 (Leve Item ID 4712)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1271.2858
$ws.Range("J39").Value = 1999.6666
$ws.Range("L39").Value = 5998.9998
$ws.Range("N39").Value = -6586.9998

# CUL row 44 (Leve Item ID 4702)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 2024.75
$ws.Range("I44").Value = 50
$ws.Range("J44").Value = 3999.5
$ws.Range("K44").Value = 150
$ws.Range("L44").Value = 11998.5
$ws.Range("M44").Value = 248
$ws.Range("N44").Value = -12794.5

# CUL row 60 (Leve Item ID 4750)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 264
$ws.Range("I60").Value = 74.85714
$ws.Range("J60").Value = 926
$ws.Range("K60").Value = 224.57142
$ws.Range("L60").Value = 2778
$ws.Range("M60").Value = 26.42858000000001
$ws.Range("N60").Value = -3280

# GSM row 49 (Leve Item ID 4232)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 33250
$ws.Range("J49").Value = 33250
$ws.Range("L49").Value = 33250
$ws.Range("N49").Value = -33618

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 24100928
$ws.Range("I132").Value = 28919348
$ws.Range("J132").Value = 8828.143
$ws.Range("K132").Value = 86758044
$ws.Range("L132").Value = 26484.429
$ws.Range("M132").Value = -86755514
$ws.Range("N132").Value = -31544.429

# LTW row 7 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3917.762
$ws.Range("I7").Value = 3624.6
$ws.Range("J7").Value = 4650.6665
$ws.Range("K7").Value = 3624.6
$ws.Range("L7").Value = 4650.6665
$ws.Range("M7").Value = -3512.6
$ws.Range("N7").Value = -4874.6665

# LTW row 16 (Leve Item ID 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 541.15
$ws.Range("J16").Value = 906.8
$ws.Range("L16").Value = 906.8
$ws.Range("N16").Value = -1246.8

# LTW row 55 (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 129.06667
$ws.Range("I55").Value = 88.09999999999999
$ws.Range("K55").Value = 88.09999999999999
$ws.Range("M55").Value = 84.90000000000001

# LTW row 126 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3917.762
$ws.Range("I126").Value = 3624.6
$ws.Range("J126").Value = 4650.6665
$ws.Range("K126").Value = 10873.8
$ws.Range("L126").Value = 13951.9995
$ws.Range("M126").Value = -8403.799999999999
$ws.Range("N126").Value = -18891.9995

# LTW row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4976027
$ws.Range("I132").Value = 5803698
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 17411094
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -17408564
$ws.Range("N132").Value = -35060

# WVR row 107 (Leve Item ID 27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1675.5186
$ws.Range("I107").Value = 944.05554
$ws.Range("J107").Value = 3138.4443
$ws.Range("K107").Value = 2832.16662
$ws.Range("L107").Value = 9415.332900000001
$ws.Range("M107").Value = -912.16662
$ws.Range("N107").Value = -13255.3329

# WVR row 113 (Leve Item ID 27752)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3338.9412
$ws.Range("I113").Value = 900.1667
$ws.Range("J113").Value = 4669.1816
$ws.Range("K113").Value = 2700.5001
$ws.Range("L113").Value = 14007.5448
$ws.Range("M113").Value = -530.5001000000002
$ws.Range("N113").Value = -18347.5448

# WVR row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3872623.8
$ws.Range("I132").Value = 4793546.5
$ws.Range("K132").Value = 14380639.5
$ws.Range("M132").Value = -14378109.5

# WVR row 133 (Leve Item ID 41869)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120
